$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 2 for the new "Posdoctorado" entry,
# pushing the existing "Programa de Comunicación..." entry down.
$ws.Rows("2:3").Insert()

# New entry: Posdoctorado
$ws.Range("A2").Value = "Posdoctorado"
$ws.Range("B2").Value = "Desde 2023"
$ws.Range("C2").Value = "Asociación Red de Mujeres Víctimas y Profesionales"
$ws.Range("D2").Value = "Bogotá, Colombia"
$ws.Range("E2").Value = "\textbf{Proyecto: } La necesidad de generar procesos de reparación social a las mujeres víctimas y sobrevivientes de violencias sexuales en el marco del conflicto armado desde el quehacer periodístico. Diversas propuestas de tratamiento según contextos"
$ws.Range("E3").Value = "Financiación del Ministerio de Ciencia Tecnología e Innovación - Minciencias"

# Update the (now shifted) "2017 - Actualmente" entry to "2017 - 2023"
$ws.Range("B4").Value = "2017 - 2023"

# Update "Actualmente" references to "2023" for the two activity rows
$ws.Range("E4").Value = "Gestión de la comunicación (4 horas semanales - 2018 - 2023)"
$ws.Range("E5").Value = "Prácticas profesionales (4 horas semanales - 2022 - 2023)"

# Adjust column widths for columns A and B
$ws.Columns("A").ColumnWidth = 43.7109375
$ws.Columns("B").ColumnWidth = 18.28515625

# Update active selection to E15
$null = $ws.Range("E15").Select()
